$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3110.875
$ws.Range("I32").Value = 2497.3333
$ws.Range("K32").Value = 2497.3333
$ws.Range("M32").Value = -2171.3333

$ws.Range("H40").Value = 4319.7
$ws.Range("I40").Value = 2999.5
$ws.Range("J40").Value = 4649.75
$ws.Range("K40").Value = 2999.5
$ws.Range("L40").Value = 4649.75
$ws.Range("M40").Value = -2824.5
$ws.Range("N40").Value = -4999.75

$ws.Range("H55").Value = 494.8
$ws.Range("I55").Value = 321.57144
$ws.Range("J55").Value = 899
$ws.Range("K55").Value = 321.57144
$ws.Range("L55").Value = 899
$ws.Range("M55").Value = -107.57144
$ws.Range("N55").Value = -1327

$ws.Range("H62").Value = 11275.409
$ws.Range("I62").Value = 14091.125
$ws.Range("J62").Value = 9666.429
$ws.Range("K62").Value = 14091.125
$ws.Range("L62").Value = 9666.429
$ws.Range("M62").Value = -13467.125
$ws.Range("N62").Value = -10914.429

$ws.Range("H65").Value = 11275.409
$ws.Range("I65").Value = 14091.125
$ws.Range("J65").Value = 9666.429
$ws.Range("K65").Value = 70455.625
$ws.Range("L65").Value = 48332.145
$ws.Range("M65").Value = -67335.625
$ws.Range("N65").Value = -54572.145

$ws.Range("H106").Value = 10794.2
$ws.Range("I106").Value = 8695.286
$ws.Range("K106").Value = 8695.286
$ws.Range("M106").Value = -8064.286

$ws.Range("H112").Value = 2420.5
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 2684.6
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 8053.799999999999
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -10269.8

$ws.Range("H116").Value = 5333.2383
$ws.Range("I116").Value = 4969.3076
$ws.Range("J116").Value = 5924.625
$ws.Range("K116").Value = 4969.3076
$ws.Range("L116").Value = 5924.625
$ws.Range("M116").Value = -1527.3076
$ws.Range("N116").Value = -12808.625

$ws.Range("H132").Value = 18616.35
$ws.Range("I132").Value = 18271.812
$ws.Range("J132").Value = 19994.5
$ws.Range("K132").Value = 54815.436
$ws.Range("L132").Value = 59983.5
$ws.Range("M132").Value = -52285.436
$ws.Range("N132").Value = -65043.5

$ws.Range("H138").Value = 3599.4055
$ws.Range("J138").Value = 6451.636
$ws.Range("L138").Value = 19354.908
$ws.Range("N138").Value = -29634.908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15544.247
$ws.Range("I32").Value = 15724.743
$ws.Range("J32").Value = 13533
$ws.Range("K32").Value = 15724.743
$ws.Range("L32").Value = 13533
$ws.Range("M32").Value = -15437.743
$ws.Range("N32").Value = -14107

$ws.Range("H61").Value = 1621.381
$ws.Range("I61").Value = 1409.375
$ws.Range("J61").Value = 2299.8
$ws.Range("K61").Value = 1409.375
$ws.Range("L61").Value = 2299.8
$ws.Range("M61").Value = -1197.375
$ws.Range("N61").Value = -2723.8

$ws.Range("H74").Value = 23789.229
$ws.Range("I74").Value = 25497.455
$ws.Range("K74").Value = 25497.455
$ws.Range("M74").Value = -24623.455

$ws.Range("H77").Value = 23789.229
$ws.Range("I77").Value = 25497.455
$ws.Range("K77").Value = 127487.275
$ws.Range("M77").Value = -123119.275

$ws.Range("H132").Value = 32342.727
$ws.Range("I132").Value = 32342.727
$ws.Range("K132").Value = 97028.181
$ws.Range("M132").Value = -94498.181

$ws.Range("H136").Value = 1621.381
$ws.Range("I136").Value = 1409.375
$ws.Range("J136").Value = 2299.8
$ws.Range("K136").Value = 4228.125
$ws.Range("L136").Value = 6899.400000000001
$ws.Range("M136").Value = -1678.125
$ws.Range("N136").Value = -11999.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 30024
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 30024
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 30024
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -30696

$ws.Range("H134").Value = 3332.6875
$ws.Range("I134").Value = 3057.3845
$ws.Range("K134").Value = 9172.1535
$ws.Range("M134").Value = -6637.1535

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3284.1428
$ws.Range("J31").Value = 3397.5
$ws.Range("L31").Value = 3397.5
$ws.Range("N31").Value = -3987.5

$ws.Range("H34").Value = 3284.1428
$ws.Range("J34").Value = 3397.5
$ws.Range("L34").Value = 3397.5
$ws.Range("N34").Value = -3801.5

$ws.Range("H58").Value = 57987.832
$ws.Range("I58").Value = 68797.53
$ws.Range("J58").Value = 3939.3333
$ws.Range("K58").Value = 68797.53
$ws.Range("L58").Value = 3939.3333
$ws.Range("M58").Value = -68594.53
$ws.Range("N58").Value = -4345.3333

$ws.Range("H86").Value = 83581.2
$ws.Range("I86").Value = 136635.33
$ws.Range("K86").Value = 136635.33
$ws.Range("M86").Value = -135512.33

$ws.Range("H89").Value = 83581.2
$ws.Range("I89").Value = 136635.33
$ws.Range("K89").Value = 683176.6499999999
$ws.Range("M89").Value = -677560.6499999999

$ws.Range("H134").Value = 42778.44
$ws.Range("I134").Value = 57742.89
$ws.Range("K134").Value = 173228.67
$ws.Range("M134").Value = -170693.67

$ws.Range("H136").Value = 57987.832
$ws.Range("I136").Value = 68797.53
$ws.Range("J136").Value = 3939.3333
$ws.Range("K136").Value = 206392.59
$ws.Range("L136").Value = 11817.9999
$ws.Range("M136").Value = -203842.59
$ws.Range("N136").Value = -16917.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 406.9091
$ws.Range("I97").Value = 194
$ws.Range("J97").Value = 528.5714
$ws.Range("K97").Value = 582
$ws.Range("L97").Value = 1585.7142
$ws.Range("M97").Value = -86
$ws.Range("N97").Value = -2577.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1760.6471
$ws.Range("I102").Value = 1294.3103
$ws.Range("K102").Value = 1294.3103
$ws.Range("M102").Value = 327.6896999999999

$ws.Range("H132").Value = 46536.78
$ws.Range("I132").Value = 53012.35
$ws.Range("J132").Value = 3366.3333
$ws.Range("K132").Value = 159037.05
$ws.Range("L132").Value = 10098.9999
$ws.Range("M132").Value = -156507.05
$ws.Range("N132").Value = -15158.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4306.1816
$ws.Range("I40").Value = 3916.8
$ws.Range("K40").Value = 3916.8
$ws.Range("M40").Value = -3780.8

$ws.Range("H61").Value = 2718.1052
$ws.Range("J61").Value = 5223.75
$ws.Range("L61").Value = 5223.75
$ws.Range("N61").Value = -5627.75

$ws.Range("H68").Value = 4298
$ws.Range("I68").Value = 4155.8
$ws.Range("K68").Value = 4155.8
$ws.Range("M68").Value = -3406.8

$ws.Range("H71").Value = 4298
$ws.Range("I71").Value = 4155.8
$ws.Range("K71").Value = 20779
$ws.Range("M71").Value = -17035

$ws.Range("H113").Value = 2718.1052
$ws.Range("J113").Value = 5223.75
$ws.Range("L113").Value = 5223.75
$ws.Range("N113").Value = -9563.75

$ws.Range("H122").Value = 3691
$ws.Range("I122").Value = 3094.9565
$ws.Range("K122").Value = 9284.869499999999
$ws.Range("M122").Value = -6834.869499999999

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()

$ws.Range("H132").Value = 66541.37
$ws.Range("I132").Value = 82786.13
$ws.Range("K132").Value = 248358.39
$ws.Range("M132").Value = -245828.39

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H129").Value = 87429
$ws.Range("J129").Value = 87429
$ws.Range("L129").Value = 87429
$ws.Range("N129").Value = -97429

$ws.Range("H132").Value = 98255.91
$ws.Range("I132").Value = 113132.055
$ws.Range("K132").Value = 339396.165
$ws.Range("M132").Value = -336866.165

$ws.Range("H136").Value = 3917.4348
$ws.Range("I136").Value = 3877.318
$ws.Range("K136").Value = 11631.954
$ws.Range("M136").Value = -9081.954000000002
